$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 3) mirroring the structure of row 2
$ws.Range("A3").Value = 42606.882835648146
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = 59
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 50
$ws.Range("F3").Value = 50
$ws.Range("G3").Value = 3572
$ws.Range("H3").Value = 3745
$ws.Range("I3").Value = 432
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 50
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1
$ws.Range("N3").Value = "Bag"
